{"js": "// Replace the date line and the 25 math-problem cells with their new values.\n// Every \"from\" string is unique across the document, so a plain\n// find-and-replace (matchCase on, no wildcards) is safe and precise.\nconst replacements = [\n  [\"2024-03-05 Tuesday\", \"2024-03-06 Wednesday\"],\n  [\"56\u00f72=28, 0\", \"91\u00f74=22, 3\"],\n  [\"42\u00f79=4, 6\", \"29\u00f77=4, 1\"],\n  [\"62\u00f75=12, 2\", \"39\u00f78=4, 7\"],\n  [\"47\u00f73=15, 2\", \"70\u00f79=7, 7\"],\n  [\"62\u00f79=6, 8\", \"22\u00f75=4, 2\"],\n  [\"92\u00f77=13, 1\", \"46\u00f76=7, 4\"],\n  [\"71\u00f73=23, 2\", \"61\u00f73=20, 1\"],\n  [\"30\u00f73=10, 0\", \"95\u00f77=13, 4\"],\n  [\"30\u00f74=7, 2\", \"64\u00f73=21, 1\"],\n  [\"16\u00f77=2, 2\", \"58\u00f79=6, 4\"],\n  [\"49\u00f78=6, 1\", \"23\u00f79=2, 5\"],\n  [\"85\u00f73=28, 1\", \"98\u00f76=16, 2\"],\n  [\"78\u00f78=9, 6\", \"87\u00f76=14, 3\"],\n  [\"27\u00f77=3, 6\", \"89\u00f75=17, 4\"],\n  [\"69\u00f72=34, 1\", \"43\u00f78=5, 3\"],\n  [\"49\u00f79=5, 4\", \"54\u00f78=6, 6\"],\n  [\"27\u00f73=9, 0\", \"98\u00f75=19, 3\"],\n  [\"38\u00f75=7, 3\", \"91\u00f75=18, 1\"],\n  [\"72\u00f78=9, 0\", \"39\u00f74=9, 3\"],\n  [\"14\u00f72=7, 0\", \"37\u00f78=4, 5\"],\n  [\"22\u00f74=5, 2\", \"31\u00f79=3, 4\"],\n  [\"38\u00f74=9, 2\", \"11\u00f72=5, 1\"],\n  [\"91\u00f78=11, 3\", \"86\u00f74=21, 2\"],\n  [\"97\u00f77=13, 6\", \"33\u00f79=3, 6\"],\n  [\"46\u00f72=23, 0\", \"46\u00f79=5, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 math-problem cells with their new values.\n# Every \"from\" string is unique across the document, so Find/Replace\n# (MatchCase on, no wildcards) targets exactly one occurrence each.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-05 Tuesday\", \"2024-03-06 Wednesday\"),\n    @(\"56\u00f72=28, 0\", \"91\u00f74=22, 3\"),\n    @(\"42\u00f79=4, 6\", \"29\u00f77=4, 1\"),\n    @(\"62\u00f75=12, 2\", \"39\u00f78=4, 7\"),\n    @(\"47\u00f73=15, 2\", \"70\u00f79=7, 7\"),\n    @(\"62\u00f79=6, 8\", \"22\u00f75=4, 2\"),\n    @(\"92\u00f77=13, 1\", \"46\u00f76=7, 4\"),\n    @(\"71\u00f73=23, 2\", \"61\u00f73=20, 1\"),\n    @(\"30\u00f73=10, 0\", \"95\u00f77=13, 4\"),\n    @(\"30\u00f74=7, 2\", \"64\u00f73=21, 1\"),\n    @(\"16\u00f77=2, 2\", \"58\u00f79=6, 4\"),\n    @(\"49\u00f78=6, 1\", \"23\u00f79=2, 5\"),\n    @(\"85\u00f73=28, 1\", \"98\u00f76=16, 2\"),\n    @(\"78\u00f78=9, 6\", \"87\u00f76=14, 3\"),\n    @(\"27\u00f77=3, 6\", \"89\u00f75=17, 4\"),\n    @(\"69\u00f72=34, 1\", \"43\u00f78=5, 3\"),\n    @(\"49\u00f79=5, 4\", \"54\u00f78=6, 6\"),\n    @(\"27\u00f73=9, 0\", \"98\u00f75=19, 3\"),\n    @(\"38\u00f75=7, 3\", \"91\u00f75=18, 1\"),\n    @(\"72\u00f78=9, 0\", \"39\u00f74=9, 3\"),\n    @(\"14\u00f72=7, 0\", \"37\u00f78=4, 5\"),\n    @(\"22\u00f74=5, 2\", \"31\u00f79=3, 4\"),\n    @(\"38\u00f74=9, 2\", \"11\u00f72=5, 1\"),\n    @(\"91\u00f78=11, 3\", \"86\u00f74=21, 2\"),\n    @(\"97\u00f77=13, 6\", \"33\u00f79=3, 6\"),\n    @(\"46\u00f72=23, 0\", \"46\u00f79=5, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #   ReplaceWith, Replace)\n    # Wrap:=1 (wdFindContinue), Replace:=2 (wdReplaceAll)\n    $found = $range.Find.Execute(\n        $oldText,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    )\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n"}
